# Add new survey response row (Hannah Minton) to Sheet1, row 17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

function Set-TextCell($ws, $row, $col, $text) {
    # Force a cell to be stored as Text even when its content looks
    # numeric/empty, without leaving a lingering custom style behind.
    if ($text -eq "") {
        # A lone apostrophe is Excel's "force text" quote-prefix; entered on
        # its own it commits as an empty Text cell.
        $ws.Cells.Item($row, $col).Value = "'"
    } else {
        # "@" (Text) number format stops a numeric-looking string like
        # "0.693" from being auto-converted to a Number on assignment.
        $ws.Cells.Item($row, $col).NumberFormat = "@"
        $ws.Cells.Item($row, $col).Value = $text
    }
    $ws.Cells.Item($row, $col).Style = "Normal"
}

$ws.Cells.Item($row, 1).Value = "Hannah Minton_20251202_125046"
Set-TextCell $ws $row 2 ""
$ws.Cells.Item($row, 3).Value = "Hannah Minton"
$ws.Cells.Item($row, 4).Value = 19
$ws.Cells.Item($row, 5).Value = "Female"
$ws.Cells.Item($row, 6).Value = "2025-12-02 12:50:47"
$ws.Cells.Item($row, 7).Value = "{
  ""portion"": 1.0,
  ""diet"": 0.5714285714285714,
  ""salt"": 0.0,
  ""fat"": 0.0,
  ""natural"": 0.2,
  ""convenience"": 0.8,
  ""price"": 1.0
}"

$ws.Cells.Item($row, 8).Value = "Nongshim Neoguri Spicy Seafood"
Set-TextCell $ws $row 9 "0.693"
$ws.Cells.Item($row, 10).Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Cells.Item($row, 11).Value = "Nissin Chow Mein Teriyaki Beef"
Set-TextCell $ws $row 12 "0.688"
$ws.Cells.Item($row, 13).Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Cells.Item($row, 14).Value = "Nongshim Shin Ramyun"
Set-TextCell $ws $row 15 "0.681"
$ws.Cells.Item($row, 16).Value = "Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio"

$ws.Cells.Item($row, 17).Value = "Amy’s Macaroni & Cheese (frozen)"
Set-TextCell $ws $row 18 "0.566"
$ws.Cells.Item($row, 19).Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

$ws.Cells.Item($row, 20).Value = "Velveeta Original Shells & Cheese (microwave cups)"
Set-TextCell $ws $row 21 "0.510"
$ws.Cells.Item($row, 22).Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"

$ws.Cells.Item($row, 23).Value = "Kraft Macaroni & Cheese Dinner"
Set-TextCell $ws $row 24 "0.440"
$ws.Cells.Item($row, 25).Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Cells.Item($row, 26).Value = "Wild Planet Wild Tuna Pasta Salad"
Set-TextCell $ws $row 27 "0.597"
$ws.Cells.Item($row, 28).Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

$ws.Cells.Item($row, 29).Value = "StarKist Chicken Creations (Chicken Salad)"
Set-TextCell $ws $row 30 "0.522"
$ws.Cells.Item($row, 31).Value = "Portátil, saludable, fácil, buena textura, sabor suave"

$ws.Cells.Item($row, 32).Value = "Jack Link’s Beef Jerky Original"
Set-TextCell $ws $row 33 "0.507"
$ws.Cells.Item($row, 34).Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# The multi-line JSON in column G makes Excel auto-expand the row height;
# AutoFit keeps it from sticking around as an explicit custom height.
$ws.Rows.Item($row).EntireRow.AutoFit()
